# Adds a new "05dec2025" daily column (F) to both the "crosstab" and
# "annot" worksheets, pushing the existing 25nov2025..30nov2025 columns
# one slot to the right (F..K -> G..L), and fills in the new day's
# counts plus the corrected 04dec2025 (column E) totals.

$wb = $excel.ActiveWorkbook

# New values for column E (04dec2025) and the newly inserted column F
# (05dec2025), keyed by data row number (row 1 is the header row).
$rowData = @{
    2  = @{ E = 0;  F = 0 }
    3  = @{ E = 21; F = 0 }
    4  = @{ E = 5;  F = 0 }
    5  = @{ E = 8;  F = 3 }
    6  = @{ E = 5;  F = 0 }
    7  = @{ E = 10; F = 5 }
    8  = @{ E = 5;  F = 0 }
    9  = @{ E = 11; F = 0 }
    10 = @{ E = 6;  F = 0 }
    11 = @{ E = 8;  F = 0 }
    12 = @{ E = 4;  F = 0 }
    13 = @{ E = 5;  F = 0 }
    14 = @{ E = 0;  F = 0 }
    15 = @{ E = 10; F = 0 }
    16 = @{ E = 11; F = 0 }
    17 = @{ E = 6;  F = 0 }
    18 = @{ E = 5;  F = 5 }
}

for ($sheetIdx = 1; $sheetIdx -le $wb.Worksheets.Count(); $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    # The "annot" sheet mirrors "crosstab" but stores every count as text,
    # leaving the cell blank instead of writing a literal "0" - match that
    # existing convention for the values we write too.
    $isAnnotSheet = ($ws.Name() -eq "annot")

    # Insert a new column in front of the old "25nov2025" column (F);
    # everything from F onward (F..K) shifts right to (G..L) along with
    # its formatting.
    $ws.Columns.Item(6).Insert()

    # Header for the freshly inserted column.
    $ws.Cells.Item(1, 6).Value = "05dec2025"

    foreach ($r in $rowData.Keys) {
        $vals = $rowData[$r]

        if ($isAnnotSheet) {
            if ($vals.E -eq 0) { $eVal = "" } else { $eVal = [string]$vals.E }
            if ($vals.F -eq 0) { $fVal = "" } else { $fVal = [string]$vals.F }
        } else {
            $eVal = $vals.E
            $fVal = $vals.F
        }

        $ws.Cells.Item($r, 5).Value = $eVal
        $ws.Cells.Item($r, 6).Value = $fVal
    }
}
